$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Insert 4 new rows before the old row 28 (old rows 28-33 shift down to 32-37) ---
$ws.Rows("27:30").Insert()

# --- Row 27 (brand new row): "modified cnn to cnn11 after Ahmed" / "Way better than cnn10!" ---
$ws.Range("B27").Value = "my laptop"
$ws.Range("C27").Value = "recognizeStressPermf"
$ws.Range("D27").Value = 11
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = $false
$ws.Range("I27").Value = "modified cnn to cnn11 after Ahmed"
$ws.Range("J27").Value = "Way better than cnn10!"
$ws.Range("K27").Value = "saved"

# J27 needs the red-font style (same style already used on sheet "16000 Bald data" J14)
$wb.Worksheets.Item(2).Range("J14").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("J27").Value = "Way better than cnn10!"

# --- Row 28 (was the "5,20" row; content replaced in place) ---
$ws.Range("B28").Value = "my laptop"
$ws.Range("C28").Value = "recognizeStressPermf"
$ws.Range("D28").Value = 11
$ws.Range("E28").Value = "1cont"
$ws.Range("F28").Value = $false
$ws.Range("G28").Value = "YES"
$ws.Range("H28").Value = "selu"
$ws.Range("H28").HorizontalAlignment = -4108
$ws.Range("I28").Value = "modified cnn to cnn11 after Ahmed"
$ws.Range("K28").ClearContents()

# J28 stays blank but picks up the same red-font style as J27
$wb.Worksheets.Item(2).Range("J14").Copy()
$ws.Range("J28").PasteSpecial(-4122)

# --- Row 29 (content replaced in place) ---
$ws.Range("B29").Value = "my laptop"
$ws.Range("C29").Value = "recognizeStressPermf"
$ws.Range("D29").Value = 11
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = $false
$ws.Range("I29").Value = "modified cnn to cnn11 after Ahmed"
$ws.Range("K29").ClearContents()

# J29 gets the plain (non-themed) font style used elsewhere (sheet "65000data" G28)
$wb.Worksheets.Item(3).Range("G28").Copy()
$ws.Range("J29").PasteSpecial(-4122)
$ws.Range("J29").Value = "worse than 11_1"

# --- Row 30 stays a genuinely blank spacer row (remove the leftover inserted cells) ---
$ws.Rows(30).Clear()

# --- Update view: selection on G32, scrolled so row 22 is at the top ---
$ws.Range("G32").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 2

"ok"
